$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.115459442138672
$ws.Range("B1").Value = 2.033017873764038
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.171510696411133
$ws.Range("E1").Value = 1.105171918869019
